$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "VALOR MORA" total (E11): 227760 -> 455520
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 455520

# ---------------------------------------------------------------------------
# 2. Update "Cant. Periodos" (F13): 1 -> 2
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = 2

# ---------------------------------------------------------------------------
# 3. Bump the "Salario Basico" column (G) for the existing 2507 rows (16-18)
#    from 1300000 to 1423500, matching row 19 which was already 1423500.
# ---------------------------------------------------------------------------
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500

# ---------------------------------------------------------------------------
# 4. Insert a second block of 4 worker rows (period 2508) right below the
#    existing block (rows 16-19). This pushes the signature block (old rows
#    24-25) down to rows 28-29 and grows the table to rows 16-23.
# ---------------------------------------------------------------------------
$ws.Rows("20:23").Insert()

# Preserve the special "last row" formatting (thicker bottom border etc.)
# that currently lives on row 19 by copying it down onto the new last row
# (23) before row 19's own formatting gets normalised below.
$ws.Range("B19:J19").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)

# Row 19 is no longer the last row of the table, so it should pick up the
# regular interior-row formatting (same as rows 16-18) instead. Rows 20-22
# (brand new, currently blank) need that same regular formatting too.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Fill in the values for the new period-2508 rows (20-23). Row 19 already
#    holds the correct (period 2507) data and was left untouched above.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73127876"
$ws.Range("D20").Value = "ARIS OMAR FERREIRA OSORIO"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73119857"
$ws.Range("D21").Value = "ORLANDO GUARDO MANJARREZ"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047450029"
$ws.Range("D22").Value = "OMAR JESUS FERREIRA ROJAS"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143409812"
$ws.Range("D23").Value = "MAIKOL MANUEL MURIEL MUÃ?OZ"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500
